# "Added Edit customer Test"
# Replace the B4 test-data value with a new hyperlinked credential
# ("Dhana@01"), shifting B5/B6 back onto their original (now
# renumbered) shared-string values, and leave the selection on B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B4 becomes the new value "Dhana@01", formatted + linked as a hyperlink.
$ws.Range("B4").Value = "Dhana@01"
$ws.Hyperlinks.Add($ws.Range("B4"), "https://example.com")

# Park the selection on B7 (matches the post-edit cursor position).
$ws.Range("B7").Select()
